$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert 3 new columns before column B (old B:E data shifts to E:H)
$ws.Range("B1:D1").EntireColumn.Insert(-4161)

# New header values for the newly inserted date columns
$ws.Range("B1").Value = "Jun_27"
$ws.Range("C1").Value = "Jun_26"
$ws.Range("D1").Value = "Jun_26"

# Fill the new columns for existing analyst rows (2-27) with "UN"
$ws.Range("B2:D27").Value = "UN"

# Add two new analyst rows with ratings for the three newest dates only
$ws.Range("A28").Value = "Benchmark"
$ws.Range("B28:D28").Value = "UN"

$ws.Range("A29").Value = "Evercore ISI"
$ws.Range("B29:D29").Value = "UN"
